$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 434, shifting existing rows 434-514 down to 435-515.
$ws.Rows("434:434").Insert()

# Populate the newly inserted row 434 with the new data point.
$ws.Cells.Item(434, 1).Value = 3
$ws.Cells.Item(434, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(434, 3).Value = "Coquimbo"
$ws.Cells.Item(434, 4).Value = 44995
$ws.Cells.Item(434, 5).Value = 5
$ws.Cells.Item(434, 6).Value = 100114013
$ws.Cells.Item(434, 7).Value = "Zanahoria"
$ws.Cells.Item(434, 8).Value = "Sin especificar"
$ws.Cells.Item(434, 9).Value = "Primera"
$ws.Cells.Item(434, 10).Value = 310
$ws.Cells.Item(434, 11).Value = 8000
$ws.Cells.Item(434, 12).Value = 8500
$ws.Cells.Item(434, 13).Value = 8242
$ws.Cells.Item(434, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(434, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(434, 16).Value = 412
$ws.Cells.Item(434, 17).Value = 20
$ws.Cells.Item(434, 18).Value = "Hortaliza"
